$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 23), continuing the subject-level stats table.
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 22
$ws.Range("C23").Value = 0.4
$ws.Range("D23").Value = 3.95
$ws.Range("E23").Value = -0.11
$ws.Range("F23").Value = 4.34
$ws.Range("G23").Value = -0.5289256198347108
$ws.Range("H23").Value = -0.48
$ws.Range("I23").Value = 1

# Column A carries the same bordered/bold/centered style as the rest of
# the index column above it (row 22) - copy formats only, not values.
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
